# Update public EPEX spot / Gaz / CO2 price workbook with the latest day.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column (BZ) for 30-aug
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Header cell - copy the style of the previous date header (BY1) so the
# new column matches the bold/centered/bordered header formatting.
$ws1.Range("BY1").Copy()
$ws1.Range("BZ1").PasteSpecial(-4122) # xlPasteFormats
$ws1.Range("BZ1").Value = "30-aug"

$ws1.Range("BZ2").Value = 67.02
$ws1.Range("BZ3").Value = 53.66
$ws1.Range("BZ4").Value = 58.64
$ws1.Range("BZ5").Value = 81.51000000000001
$ws1.Range("BZ6").Value = 37.79
$ws1.Range("BZ7").Value = 25.32
$ws1.Range("BZ8").Value = 63.72
$ws1.Range("BZ9").Value = 54.56
$ws1.Range("BZ10").Value = 70.56999999999999
$ws1.Range("BZ11").Value = 32
$ws1.Range("BZ12").Value = 12.82
$ws1.Range("BZ13").Value = 5.23
$ws1.Range("BZ14").Value = 0
$ws1.Range("BZ15").Value = -0.01
$ws1.Range("BZ16").Value = -0.01
$ws1.Range("BZ17").Value = 0
$ws1.Range("BZ18").Value = 2.54
$ws1.Range("BZ19").Value = 7.1
$ws1.Range("BZ20").Value = 18.99
$ws1.Range("BZ21").Value = 57.97
$ws1.Range("BZ22").Value = 69.68000000000001
$ws1.Range("BZ23").Value = 49.08
$ws1.Range("BZ24").Value = 45.08
$ws1.Range("BZ25").Value = 40.38

# ------------------------------------------------------------------
# Sheet "Gaz": append a new row (75) for 2025-08-28
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the date column to be stored as plain text (matches the rest of
# column A, which holds "yyyy-mm-dd" strings, not real Excel dates),
# then drop the number-format override so the cell keeps the sheet's
# default (unstyled) look, same as every other row in the column.
$ws2.Range("A75").NumberFormat = "@"
$ws2.Range("A75").Value = "2025-08-28"
$ws2.Range("A75").ClearFormats()
$ws2.Range("B75").Value = 30.55

# ------------------------------------------------------------------
# Sheet "CO2": append a new row (75) for 2025-08-28
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A75").NumberFormat = "@"
$ws3.Range("A75").Value = "2025-08-28"
$ws3.Range("A75").ClearFormats()
$ws3.Range("B75").Value = 70.95
